# Edit: rename headers and add a "PO Forecast" sheet with forecast data.
$wb = $excel.ActiveWorkbook

# --- 1) Rename "Requested quantity" headers on the existing sheets ---
$weekly = $wb.Worksheets.Item("Weekly Quantity")
$weekly.Range("B1").Value = "Weekly_PO_Qty"

$monthly = $wb.Worksheets.Item("Monthly Trend")
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet after the last existing sheet ---
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$forecast = $wb.Worksheets.Add($null, $lastSheet)
$forecast.Name = "PO Forecast"

# --- 3) Header row ---
$forecast.Range("A1").Value = "ds"
$forecast.Range("B1").Value = "PO_Forecast"
$forecast.Range("C1").Value = "yhat_lower"
$forecast.Range("D1").Value = "yhat_upper"

# --- 4) Forecast data rows (A2:D37) ---
$data = New-Object 'object[,]' 36,4
$data[0,0] = 45319.99999999999; $data[0,1] = 64; $data[0,2] = -644.8464440196043; $data[0,3] = 798.2887908857468
$data[1,0] = 45333.99999999999; $data[1,1] = 105; $data[1,2] = -548.4108590788097; $data[1,3] = 836.6306896970494
$data[2,0] = 45340.99999999999; $data[2,1] = 126; $data[2,2] = -588.8944470351353; $data[2,3] = 812.307520138449
$data[3,0] = 45354.99999999999; $data[3,1] = 167; $data[3,2] = -510.0683514502867; $data[3,3] = 864.2883804374235
$data[4,0] = 45361.99999999999; $data[4,1] = 188; $data[4,2] = -474.936317688177; $data[4,3] = 856.7286697409528
$data[5,0] = 45368.99999999999; $data[5,1] = 208; $data[5,2] = -461.525415591232; $data[5,3] = 898.5853937791901
$data[6,0] = 45375.99999999999; $data[6,1] = 229; $data[6,2] = -445.005412811373; $data[6,3] = 881.667935894822
$data[7,0] = 45417.99999999999; $data[7,1] = 352; $data[7,2] = -291.2590316997014; $data[7,3] = 1074.924800920107
$data[8,0] = 45431.99999999999; $data[8,1] = 393; $data[8,2] = -276.1856648451276; $data[8,3] = 1103.93157259332
$data[9,0] = 45438.99999999999; $data[9,1] = 414; $data[9,2] = -260.5592849257087; $data[9,3] = 1117.707258705977
$data[10,0] = 45445.99999999999; $data[10,1] = 434; $data[10,2] = -237.7544814948554; $data[10,3] = 1137.588245278144
$data[11,0] = 45459.99999999999; $data[11,1] = 475; $data[11,2] = -186.1918071625985; $data[11,3] = 1105.306287322947
$data[12,0] = 45466.99999999999; $data[12,1] = 496; $data[12,2] = -190.3278084714825; $data[12,3] = 1154.918811252843
$data[13,0] = 45473.99999999999; $data[13,1] = 517; $data[13,2] = -216.4804036609116; $data[13,3] = 1189.517707385586
$data[14,0] = 45480.99999999999; $data[14,1] = 537; $data[14,2] = -158.1082561236182; $data[14,3] = 1213.853311806576
$data[15,0] = 45487.99999999999; $data[15,1] = 558; $data[15,2] = -78.91383973806002; $data[15,3] = 1282.033045242702
$data[16,0] = 45494.99999999999; $data[16,1] = 578; $data[16,2] = -69.8271760656192; $data[16,3] = 1201.857415658654
$data[17,0] = 45501.99999999999; $data[17,1] = 599; $data[17,2] = -97.47732997127233; $data[17,3] = 1270.621683277453
$data[18,0] = 45508.99999999999; $data[18,1] = 619; $data[18,2] = -71.85794322542269; $data[18,3] = 1273.025405627073
$data[19,0] = 45529.99999999999; $data[19,1] = 681; $data[19,2] = -30.2438870268317; $data[19,3] = 1406.331790503366
$data[20,0] = 45536.99999999999; $data[20,1] = 702; $data[20,2] = 52.90563972202438; $data[20,3] = 1397.053898542548
$data[21,0] = 45543.99999999999; $data[21,1] = 722; $data[21,2] = 16.2672424605577; $data[21,3] = 1396.192909763779
$data[22,0] = 45550.99999999999; $data[22,1] = 743; $data[22,2] = 64.97535903943036; $data[22,3] = 1413.423378857597
$data[23,0] = 45564.99999999999; $data[23,1] = 784; $data[23,2] = 83.71288839213392; $data[23,3] = 1481.730694538099
$data[24,0] = 45571.99999999999; $data[24,1] = 805; $data[24,2] = 71.99045920885244; $data[24,3] = 1456.276497378261
$data[25,0] = 45592.99999999999; $data[25,1] = 866; $data[25,2] = 103.2102338275231; $data[25,3] = 1551.689865190599
$data[26,0] = 45599.99999999999; $data[26,1] = 887; $data[26,2] = 258.6410408341284; $data[26,3] = 1547.256331541648
$data[27,0] = 45613.99999999999; $data[27,1] = 928; $data[27,2] = 198.4307432615962; $data[27,3] = 1620.932958302436
$data[28,0] = 45620.99999999999; $data[28,1] = 949; $data[28,2] = 263.5166481595454; $data[28,3] = 1594.603342233498
$data[29,0] = 45627.99999999999; $data[29,1] = 969; $data[29,2] = 230.6846825930001; $data[29,3] = 1630.024111042661
$data[30,0] = 45634.99999999999; $data[30,1] = 990; $data[30,2] = 307.5101290241762; $data[30,3] = 1653.474379929041
$data[31,0] = 45641.99999999999; $data[31,1] = 1010; $data[31,2] = 350.1838816279779; $data[31,3] = 1721.7444893009
$data[32,0] = 45648.99999999999; $data[32,1] = 1031; $data[32,2] = 303.3474059943389; $data[32,3] = 1706.854466069867
$data[33,0] = 45655.99999999999; $data[33,1] = 1051; $data[33,2] = 391.488078088915; $data[33,3] = 1744.981548179792
$data[34,0] = 45662.99999999999; $data[34,1] = 1072; $data[34,2] = 374.1903571163962; $data[34,3] = 1762.532757665545
$data[35,0] = 45669.99999999999; $data[35,1] = 1093; $data[35,2] = 435.4905147058968; $data[35,3] = 1757.180887716079
$forecast.Range("A2:D37").Value = $data

# --- 5) Match formatting of the source sheets: bold/centered header style,
#        and the date number format on column A ---
$weekly.Range("A1:B1").Copy()
$forecast.Range("A1:D1").PasteSpecial(-4122)

$weekly.Range("A2").Copy()
$forecast.Range("A2:A37").PasteSpecial(-4122)

$excel.CutCopyMode = 0
